$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Convolution:" heading - paragraph that used to hold just a
#    bold space becomes the bold heading "Convolution:"
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(77).Range
$p1.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "Convolution:", 2)

# ------------------------------------------------------------------
# 2) Paragraph describing the length of the convolution result.
#    Formatting goes from bold to regular, and new text is added.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(78).Range
$p2.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "The length of the convolution result is the length of the signal plus the length of the convolution kernel minus one. ", 2)
$d.Paragraphs(78).Range.Bold = 0

# ------------------------------------------------------------------
# 3) Paragraph about flipping the kernel before convolving.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(79).Range
$p3.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "When doing a convolution, we should flip the kernel before applying.it is flipped because convolution theorem says so but also, when you flip, then the convolution with an impulse response function of a system gives you the response of that system. If you don't flip, the response comes out backwards.", 2)
$d.Paragraphs(79).Range.Bold = 0

# ------------------------------------------------------------------
# 4) New paragraph: convolution theorem (bold), inserted right after
#    paragraph 79.
# ------------------------------------------------------------------
$d.Paragraphs(79).Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(80)
$newPara.Range.InsertAfter("Instead of convolution in the time domain we can transform both signal and kernel to the frequency domain and multiply them. This is the convolution theorem.")
$d.Paragraphs(80).Range.Bold = -1

# ------------------------------------------------------------------
# 5) Existing (now shifted) empty bold paragraph gets more bold text
#    appended after its leading space.
# ------------------------------------------------------------------
$d.Paragraphs(81).Range.InsertAfter("Convolution with time domain gaussian gives a smoothing filter. Convolution with a frequency domain gaussian gives a narrowband filter.")
$d.Paragraphs(81).Range.Bold = -1
